# daily auto push: 2026-01-25 09:32 UTC
# Insert a new daily-stats row for 2026/01/25 (日) right before the
# 2026/12/29 block, shifting every following row down by one, and
# extending the sheet from A1:D736 to A1:D737.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 695..736 down to 696..737, leaving a blank row at 695.
$ws.Rows.Item(695).Insert()

# Fill in the new row with the new day's data. Force column A to be
# stored as literal text (matches the rest of the "date" column, which
# is inline/shared text, not a real date), then restore the cell to an
# unstyled "Normal" state so no stray formatting is left behind.
$newRow = $ws.Range("A695")
$newRow.NumberFormat = "@"
$newRow.Value = "2026/01/25"
$newRow.Style = "Normal"

$ws.Range("B695").Value = "日"
$ws.Range("C695").Value = 16
$ws.Range("D695").Value = 18
